$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2:D3").Value = -0.269
$ws.Range("E2:E3").Value = -0.195
$ws.Range("G2:G3").Value = -2.292134831460674
$ws.Range("H2:H3").Value = -2.292134831460674
$ws.Range("I2:I3").Value = -3.707865168539326
$ws.Range("J2:J3").Value = -3.707865168539326
$ws.Range("K2:K3").Value = 1.26
$ws.Range("L2:L3").Value = 1.415730337078652
$ws.Range("M2:M3").Value = 1.61
$ws.Range("N2:N3").Value = 0.02336719883889695
$ws.Range("O2:O3").Value = 1.277777777777778
$ws.Range("P2:P3").Value = 1.61
$ws.Range("Q2:Q3").Value = 0.02336719883889695
$ws.Range("R2:R3").Value = 1.277777777777778
$ws.Range("U2:U3").Value = 6.24
$ws.Range("V2:V3").Value = 0.09056603773584905
$ws.Range("W2:W3").Value = 0.03239074550128535
$ws.Range("X2:X3").Value = 0.06782803348748273
$ws.Range("Y2:Y3").Value = -0.03543728798619739
$ws.Range("Z2:Z3").Value = 0.02470163752428532
$ws.Range("AA2:AA3").Value = -0.09159034138218151
$ws.Range("AB2:AB3").Value = 0.06534726665539836
$ws.Range("AC2:AC3").Value = -0.1569376080375799
$ws.Range("AD2:AD3").Value = 4.58
$ws.Range("AF2:AF3").Value = 4.58
$ws.Range("AG2:AG3").Value = -1.66
$ws.Range("AH2:AH3").Value = 0.06232988568317909
$ws.Range("AI2:AI3").Value = 0.108325449385052
$ws.Range("AJ2:AJ3").Value = -0.02468768590124926
$ws.Range("AK2:AK3").Value = -0.04605993340732519
$ws.Range("AL2:AL3").Value = 0.378
$ws.Range("AM2:AM3").Value = 0.11
$ws.Range("AN2:AN3").Value = -2.301507537688442
$ws.Range("AO2:AO3").Value = -8.730158730158729
$ws.Range("AP2:AP3").Value = 0.8341708542713568
$ws.Range("AQ2:AQ3").Value = -30
